$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1426048780487805
$ws.Range("V2").Value = 0.0002220611491829204
$ws.Range("Z2").Value = -0.1980364780318399
$ws.Range("AB2").Value = -891.8105610122262
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -891.8105610122262

# Row 3
$ws.Range("T3").Value = 0.1423414634146342
$ws.Range("V3").Value = 0.0002565454225970831
$ws.Range("Z3").Value = -0.1349408809612664
$ws.Range("AB3").Value = -525.9921599661418
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -525.9921599661418

# Row 4
$ws.Range("T4").Value = 0.1477268292682927
$ws.Range("V4").Value = 0.0002491214197856265
$ws.Range("Z4").Value = -0.2051903921641787
$ws.Range("AB4").Value = -823.6561606816015
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -823.6561606816015

# Row 5
$ws.Range("T5").Value = 0.1470341463414634
$ws.Range("V5").Value = 0.0002082674398172554
$ws.Range("Z5").Value = -0.2400119381765305
$ws.Range("AB5").Value = -1152.421801444957
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -1152.421801444957

# Row 6
$ws.Range("T6").Value = 0.1477560975609756
$ws.Range("V6").Value = 0.0001607142857142857
$ws.Range("Z6").Value = -0.1997284282732178
$ws.Range("AB6").Value = -1242.754664811133
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1242.754664811133

# Row 7
$ws.Range("T7").Value = 0.1455317073170732
$ws.Range("V7").Value = 0.0002247627833421192
$ws.Range("Z7").Value = -0.1011027677634835
$ws.Range("AB7").Value = -449.8198779180959
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -449.8198779180959

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = -0.002924269337526327
$ws.Range("AB8").Value = "-Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "-Inf"

# Row 9
$ws.Range("T9").Value = 0.1426048780487805
$ws.Range("V9").Value = 0.0002220611491829204
$ws.Range("Z9").Value = 0.1567247502110327
$ws.Range("AB9").Value = 705.7729404162113
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 705.7729404162113

# Row 10
$ws.Range("T10").Value = 0.1423414634146342
$ws.Range("V10").Value = 0.0002565454225970831
$ws.Range("Z10").Value = 0.1154633257018739
$ws.Range("AB10").Value = 450.0697168283318
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 450.0697168283318

# Row 11
$ws.Range("T11").Value = 0.1477268292682927
$ws.Range("V11").Value = 0.0002491214197856265
$ws.Range("Z11").Value = 0.2054295288385532
$ws.Range("AB11").Value = 824.6160808465569
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 824.6160808465569

# Row 12
$ws.Range("T12").Value = 0.1470341463414634
$ws.Range("V12").Value = 0.0002082674398172554
$ws.Range("Z12").Value = 0.2431869543741488
$ws.Range("AB12").Value = 1167.666700985683
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 1167.666700985683

# Row 13
$ws.Range("T13").Value = 0.1477560975609756
$ws.Range("V13").Value = 0.0001607142857142857
$ws.Range("Z13").Value = 0.180273597705175
$ws.Range("AB13").Value = 1121.702385721089
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 1121.702385721089

# Row 14
$ws.Range("T14").Value = 0.1455317073170732
$ws.Range("V14").Value = 0.0002247627833421192
$ws.Range("Z14").Value = 0.1054619262010817
$ws.Range("AB14").Value = 469.2143629515145
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 469.2143629515145

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = 0.00008408941819272763
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"

